$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values (B2:E2)
$ws.Range("B2").Value = 14.336164359257241
$ws.Range("C2").Value = 31.099057571451681
$ws.Range("D2").Value = 14.999933945741793
$ws.Range("E2").Value = 27.0788659403709

# Update row 3 values (B3:E3)
$ws.Range("B3").Value = 15.606197544409987
$ws.Range("C3").Value = 19.790493933408463
$ws.Range("D3").Value = 13.824821612113388
$ws.Range("E3").Value = 16.947694771919032

# Update the sheet selection to match the new active range
[void]$ws.Range("B1:E3").Select()
